$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

function Set-CellText($table, $row, $col, $newText) {
    $cell = $table.Cell($row, $col)
    $cell.Range.Text = $newText
}

# Row 5: Bevacizumab cost row
Set-CellText $t 5 2 "1517.69"
Set-CellText $t 5 3 "1214.15"
Set-CellText $t 5 4 "1821.23"
Set-CellText $t 5 6 "GAMMA(100.00, 15.18)"

# Row 24: Leukopenia adverse event disutility row
Set-CellText $t 24 2 "0.45"
Set-CellText $t 24 3 "0.36"
Set-CellText $t 24 4 "0.54"
Set-CellText $t 24 6 "alpha_u_AE1, beta_u_AE1 (54.55, 66.67)"

# Row 25: Diarrhea adverse event disutility row
Set-CellText $t 25 2 "0.19"
Set-CellText $t 25 3 "0.15"
Set-CellText $t 25 4 "0.23"
Set-CellText $t 25 6 "alpha_u_AE2, beta_u_AE2 (80.81, 344.51)"

# Row 26: Vomiting adverse event disutility row
Set-CellText $t 26 2 "0.36"
Set-CellText $t 26 3 "0.29"
Set-CellText $t 26 4 "0.43"
Set-CellText $t 26 6 "alpha_u_AE3, beta_u_AE3 (63.64, 113.14)"

# Row 29: hazard ratio row label rename
Set-CellText $t 29 1 "PFS to Dead under the Experimental Strategy"

# Narrow the last table column (tblGrid gridCol 4741 -> 4618 dxa = 230.9pt)
$t.Columns.Item(6).Width = 230.9
